$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A10").Value = "460A14400"
$ws.Range("B10").Value = "Stellingen"
$ws.Range("C10").Value = "14a798eb-7613-4fab-a588-7b642203a559"

$ws.Range("A10").HorizontalAlignment = $ws.Range("A9").HorizontalAlignment

$ws.Range("B6").Select()
